# Updates cryptos list values (Price and Volume(1h) columns)
# generated from OOXML diff - preserves text cell type & original styling
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.870.34"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = "'3.368.24"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'569.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').Value = "'138.88"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').Value = "'7.66"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('E10').Value = '  -2.38%  '
$ws.Range('D11').Value = "'0.382"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.86%  '
$ws.Range('D12').Value = "'3.942.94"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').Value = "'27.74"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').Value = "'3.367.32"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = "'0.0000168"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.00%  '
$ws.Range('D17').Value = "'60.944.02"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = "'6.07"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.92%  '
$ws.Range('E19').Value = '  -3.66%  '
$ws.Range('D20').Value = "'8.87"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.73%  '
$ws.Range('D21').Value = "'381.87"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('D22').Value = "'75.41"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  -5.61%  '
$ws.Range('D26').Value = "'0.188"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.35%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = "'7.15"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.82%  '
$ws.Range('D29').Value = "'7.80"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.70%  '
$ws.Range('E30').Value = '  -2.05%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = "'1.34"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.24%  '
$ws.Range('D33').Value = "'22.95"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.16%  '
$ws.Range('D34').Value = "'6.85"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('D35').Value = "'167.09"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('D36').Value = "'4.92"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.95%  '
$ws.Range('D37').Value = "'3.403.62"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('E38').Value = '  -3.92%  '
$ws.Range('D39').Value = "'0.0761"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.34%  '
$ws.Range('D40').Value = "'25.30"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.59%  '
$ws.Range('E41').Value = '  -1.11%  '
$ws.Range('E42').Value = '  -2.70%  '
$ws.Range('E43').Value = '  -3.70%  '
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('D45').Value = "'2.449.74"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.29%  '
$ws.Range('E47').Value = '  -3.97%  '
$ws.Range('D48').Value = "'22.17"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.89%  '
$ws.Range('D49').Value = "'0.0257"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.86%  '
$ws.Range('E50').Value = '  -2.54%  '
$ws.Range('E51').Value = '  -3.99%  '
